$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-85 down to 84-86
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new Valencia record
$ws.Range("A83").Value = 1
$ws.Range("B83").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C83").Value = "Arica y Parinacota"
$ws.Range("D83").NumberFormat = $ws.Range("D84").NumberFormat
$ws.Range("D83").Value = 44714
$ws.Range("E83").Value = 15
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100102
$ws.Range("H83").Value = "Cítricos"
$ws.Range("I83").Value = 100102005
$ws.Range("J83").Value = "Naranja"
$ws.Range("K83").Value = "Valencia"
$ws.Range("L83").Value = "Segunda"
$ws.Range("M83").Value = 270
$ws.Range("N83").Value = 850
$ws.Range("O83").Value = 900
$ws.Range("P83").Value = 875
$ws.Range("Q83").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R83").Value = "Región de Coquimbo"
$ws.Range("S83").Value = 875
$ws.Range("T83").Value = 1
